$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 300
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 300
$ws.Range("M18").Value = -16
$ws.Range("N18").Value = -868
$ws.Range("H40").Value = 4142.4116
$ws.Range("I40").Value = 2936.4167
$ws.Range("J40").Value = 7036.8
$ws.Range("K40").Value = 2936.4167
$ws.Range("L40").Value = 7036.8
$ws.Range("M40").Value = -2761.4167
$ws.Range("N40").Value = -7386.8
$ws.Range("H55").Value = 738.2941
$ws.Range("I55").Value = 450
$ws.Range("J55").Value = 1430.2
$ws.Range("K55").Value = 450
$ws.Range("L55").Value = 1430.2
$ws.Range("M55").Value = -236
$ws.Range("N55").Value = -1858.2
$ws.Range("H86").Value = 3985.6428
$ws.Range("J86").Value = 3916.5833
$ws.Range("L86").Value = 3916.5833
$ws.Range("N86").Value = -6162.5833
$ws.Range("H89").Value = 3985.6428
$ws.Range("J89").Value = 3916.5833
$ws.Range("L89").Value = 19582.9165
$ws.Range("N89").Value = -30814.9165
$ws.Range("H98").Value = 1697.7333
$ws.Range("I98").Value = 1787
$ws.Range("J98").Value = 1519.2
$ws.Range("K98").Value = 1787
$ws.Range("L98").Value = 1519.2
$ws.Range("M98").Value = -289
$ws.Range("N98").Value = -4515.2
$ws.Range("H113").Value = 8252.733
$ws.Range("I113").Value = 7100
$ws.Range("K113").Value = 7100
$ws.Range("M113").Value = -3846
$ws.Range("H120").Value = 52615
$ws.Range("J120").Value = 52615
$ws.Range("L120").Value = 52615
$ws.Range("N120").Value = -62291
$ws.Range("H122").Value = 1697.7333
$ws.Range("I122").Value = 1787
$ws.Range("J122").Value = 1519.2
$ws.Range("K122").Value = 5361
$ws.Range("L122").Value = 4557.6
$ws.Range("M122").Value = -2911
$ws.Range("N122").Value = -9457.6
$ws.Range("H137").Value = 25643512
$ws.Range("I137").Value = 37038960
$ws.Range("J137").Value = 3752.9167
$ws.Range("K137").Value = 111116880
$ws.Range("L137").Value = 11258.7501
$ws.Range("M137").Value = -111114330
$ws.Range("N137").Value = -16358.7501

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15628418
$ws.Range("I32").Value = 17546750
$ws.Range("J32").Value = 7715.5713
$ws.Range("K32").Value = 17546750
$ws.Range("L32").Value = 7715.5713
$ws.Range("M32").Value = -17546463
$ws.Range("N32").Value = -8289.5713
$ws.Range("H113").Value = 62395.5
$ws.Range("J113").Value = 62395.5
$ws.Range("L113").Value = 62395.5
$ws.Range("N113").Value = -71073.5
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H119").Value = 65154
$ws.Range("J119").Value = 65154
$ws.Range("L119").Value = 65154
$ws.Range("N119").Value = -74830
$ws.Range("H122").Value = 1781.3529
$ws.Range("I122").Value = 1810.7241
$ws.Range("K122").Value = 5432.1723
$ws.Range("M122").Value = -2982.1723
$ws.Range("H132").Value = 6581.8335
$ws.Range("I132").Value = 2420.5908
$ws.Range("K132").Value = 7261.7724
$ws.Range("M132").Value = -4731.7724
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4668.8335
$ws.Range("I86").Value = 3509.7856
$ws.Range("K86").Value = 3509.7856
$ws.Range("M86").Value = -2386.7856
$ws.Range("H89").Value = 4668.8335
$ws.Range("I89").Value = 3509.7856
$ws.Range("K89").Value = 17548.928
$ws.Range("M89").Value = -11932.928
$ws.Range("H94").Value = 1270.0555
$ws.Range("I94").Value = 798.2
$ws.Range("K94").Value = 798.2
$ws.Range("M94").Value = -347.2
$ws.Range("H107").Value = 1635.0869
$ws.Range("I107").Value = 1391.2273
$ws.Range("K107").Value = 1391.2273
$ws.Range("M107").Value = 528.7727
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H29").Value = 9998.5
$ws.Range("J29").Value = 9998.5
$ws.Range("L29").Value = 9998.5
$ws.Range("N29").Value = -10584.5
$ws.Range("H31").Value = 41449.258
$ws.Range("J31").Value = 64949.47
$ws.Range("L31").Value = 64949.47
$ws.Range("N31").Value = -65539.47
$ws.Range("H34").Value = 41449.258
$ws.Range("J34").Value = 64949.47
$ws.Range("L34").Value = 64949.47
$ws.Range("N34").Value = -65353.47
$ws.Range("H110").Value = 85694
$ws.Range("J110").Value = 98750
$ws.Range("L110").Value = 98750
$ws.Range("N110").Value = -106930
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("H132").Value = 4864.143
$ws.Range("I132").Value = 2187.973
$ws.Range("J132").Value = 24667.8
$ws.Range("K132").Value = 6563.919
$ws.Range("L132").Value = 74003.39999999999
$ws.Range("M132").Value = -4033.919
$ws.Range("N132").Value = -79063.39999999999

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 23.625
$ws.Range("I12").Value = 101
$ws.Range("K12").Value = 303
$ws.Range("M12").Value = -130
$ws.Range("H140").Value = 3563.7144
$ws.Range("I140").Value = 1789.8
$ws.Range("K140").Value = 5369.4
$ws.Range("M140").Value = -189.3999999999996

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10876.65
$ws.Range("I70").Value = 8338.833000000001
$ws.Range("J70").Value = 14683.375
$ws.Range("K70").Value = 8338.833000000001
$ws.Range("L70").Value = 14683.375
$ws.Range("M70").Value = -8068.833000000001
$ws.Range("N70").Value = -15223.375
$ws.Range("H73").Value = 10876.65
$ws.Range("I73").Value = 8338.833000000001
$ws.Range("J73").Value = 14683.375
$ws.Range("K73").Value = 8338.833000000001
$ws.Range("L73").Value = 14683.375
$ws.Range("M73").Value = -7402.833000000001
$ws.Range("N73").Value = -16555.375
$ws.Range("H95").Value = 0
$ws.Range("J95").Value = 0
$ws.Range("L95").Value = 0
$ws.Range("N95").ClearContents()
$ws.Range("H114").Value = 94000
$ws.Range("J114").Value = 94000
$ws.Range("L114").Value = 94000
$ws.Range("N114").Value = -102678
$ws.Range("H120").Value = 52631
$ws.Range("J120").Value = 52631
$ws.Range("L120").Value = 52631
$ws.Range("N120").Value = -62307
$ws.Range("H126").Value = 4130.1577
$ws.Range("I126").Value = 2870.6667
$ws.Range("J126").Value = 6289.2856
$ws.Range("K126").Value = 8612.000100000001
$ws.Range("L126").Value = 18867.8568
$ws.Range("M126").Value = -6142.000100000001
$ws.Range("N126").Value = -23807.8568
$ws.Range("H132").Value = 237114.1
$ws.Range("I132").Value = 270236.28
$ws.Range("J132").Value = 5258.7144
$ws.Range("K132").Value = 810708.8400000001
$ws.Range("L132").Value = 15776.1432
$ws.Range("M132").Value = -808178.8400000001
$ws.Range("N132").Value = -20836.1432

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 7386.759
$ws.Range("I7").Value = 4591.8184
$ws.Range("J7").Value = 16170.857
$ws.Range("K7").Value = 4591.8184
$ws.Range("L7").Value = 16170.857
$ws.Range("M7").Value = -4479.8184
$ws.Range("N7").Value = -16394.857
$ws.Range("H122").Value = 6092.9287
$ws.Range("I122").Value = 3421.889
$ws.Range("J122").Value = 10900.8
$ws.Range("K122").Value = 10265.667
$ws.Range("L122").Value = 32702.4
$ws.Range("M122").Value = -7815.667000000001
$ws.Range("N122").Value = -37602.39999999999
$ws.Range("H126").Value = 7386.759
$ws.Range("I126").Value = 4591.8184
$ws.Range("J126").Value = 16170.857
$ws.Range("K126").Value = 13775.4552
$ws.Range("L126").Value = 48512.571
$ws.Range("M126").Value = -11305.4552
$ws.Range("N126").Value = -53452.571
$ws.Range("H132").Value = 3859.2778
$ws.Range("I132").Value = 2070.818
$ws.Range("K132").Value = 6212.454000000001
$ws.Range("M132").Value = -3682.454000000001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4486.6924
$ws.Range("I81").Value = 3277.0833
$ws.Range("J81").Value = 19002
$ws.Range("K81").Value = 6554.1666
$ws.Range("L81").Value = 38004
$ws.Range("M81").Value = -5493.1666
$ws.Range("N81").Value = -40126
$ws.Range("H84").Value = 4486.6924
$ws.Range("I84").Value = 3277.0833
$ws.Range("J84").Value = 19002
$ws.Range("K84").Value = 32770.833
$ws.Range("L84").Value = 190020
$ws.Range("M84").Value = -27466.833
$ws.Range("N84").Value = -200628
$ws.Range("H113").Value = 481.5238
$ws.Range("I113").Value = 221
$ws.Range("K113").Value = 663
$ws.Range("M113").Value = 1507
$ws.Range("H132").Value = 7789.032
$ws.Range("I132").Value = 3395
$ws.Range("K132").Value = 10185
$ws.Range("M132").Value = -7655
$ws.Range("H133").Value = 49500
$ws.Range("J133").Value = 49500
$ws.Range("L133").Value = 49500
$ws.Range("N133").Value = -59620
